$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.375.31"
$ws.Range("E2").Value = "  -1.11%  "
$ws.Range("D3").Value = "2.987.14"
$ws.Range("E3").Value = "  -0.65%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.65%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").Value = "2.984.96"
$ws.Range("E8").Value = "  -0.71%  "
$ws.Range("B9").Value = "XRP"
$ws.Range("C9").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.514"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.70%  "
$ws.Range("E10").Value = "  -1.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.06"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.09%  "
$ws.Range("E12").Value = "  +2.90%  "
$ws.Range("E13").Value = "  -0.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.06%  "
$ws.Range("E15").Value = "  +2.80%  "
$ws.Range("D16").Value = "3.482.04"
$ws.Range("E16").Value = "  -0.55%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.90"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.97%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "61.393.88"
$ws.Range("E18").Value = "  -1.13%  "
$ws.Range("D19").Value = "2.987.52"
$ws.Range("E19").Value = "  -0.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "447.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.05%  "
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.23%  "
$ws.Range("E26").Value = "  -2.52%  "
$ws.Range("E27").Value = "  -1.77%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("E29").Value = "  +2.25%  "
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.18"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.63%  "
$ws.Range("E32").Value = "  -1.35%  "
$ws.Range("E33").Value = "  -0.34%  "
$ws.Range("E34").Value = "  +1.61%  "
$ws.Range("D35").Value = "0.0₃0818"
$ws.Range("E35").Value = "  +4.03%  "
$ws.Range("E36").Value = "  -0.49%  "
$ws.Range("E37").Value = "  +0.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "50.17"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.04"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.03"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.85"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.120"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "386.45"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.41%  "
$ws.Range("E44").Value = "  -0.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.268"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.07%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "2.692.53"
$ws.Range("E46").Value = "  -1.86%  "
$ws.Range("B47").Value = "Arweave"
$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "38.04"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "130.38"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.43%  "
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("E50").Value = "  -0.70%  "
$ws.Range("E51").Value = "  -0.32%  "
